function RGB($r, $g, $b) { return $r + ($g * 256) + ($b * 65536) }

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Row 5 / Row 6 rework: the old "buffer" row 6 (D6:E6, red fill)
# is folded up into row 5, and row 5's H:I cells get a brand new
# accent colour (instead of the old blue).
# ---------------------------------------------------------------
$ws.Range("D5:E5").Interior.Color = (RGB 255 0 0)
$ws.Range("H5:I5").Interior.ThemeColor = 10
$ws.Range("D6:E6").Clear() | Out-Null

# ---------------------------------------------------------------
# Rows 9-15 keep their text/colour but become two-column
# (D:E) centered, merged blocks. Row 16 is a brand new entry.
# ---------------------------------------------------------------
$ws.Range("D9:E9").Merge() | Out-Null
$ws.Range("D9:E9").HorizontalAlignment = -4108
$ws.Range("D9").Interior.ThemeColor = 6

$ws.Range("D10:E10").Merge() | Out-Null
$ws.Range("D10:E10").HorizontalAlignment = -4108
$ws.Range("D10").Interior.ThemeColor = 5

$ws.Range("D11:E11").Merge() | Out-Null
$ws.Range("D11:E11").HorizontalAlignment = -4108
$ws.Range("D11").Interior.Color = (RGB 146 208 80)

$ws.Range("D12:E12").Merge() | Out-Null
$ws.Range("D12:E12").HorizontalAlignment = -4108
$ws.Range("D12").Interior.Color = (RGB 255 255 0)

$ws.Range("D13").Value = "Testfälle"
$ws.Range("D13:E13").Merge() | Out-Null
$ws.Range("D13:E13").HorizontalAlignment = -4108
$ws.Range("D13").Interior.Color = (RGB 255 0 0)

$ws.Range("D14:E14").Merge() | Out-Null
$ws.Range("D14:E14").HorizontalAlignment = -4108
$ws.Range("D14").Interior.Color = (RGB 0 32 96)

$ws.Range("D15:E15").Merge() | Out-Null
$ws.Range("D15:E15").HorizontalAlignment = -4108
$ws.Range("E2").Copy() | Out-Null
$ws.Range("D15").PasteSpecial(-4122) | Out-Null

$ws.Range("D16").Value = "Testprotokoll"
$ws.Range("D16:E16").Merge() | Out-Null
$ws.Range("D16:E16").HorizontalAlignment = -4108
$ws.Range("D16").Interior.ThemeColor = 10

# ---------------------------------------------------------------
# Misc view bits
# ---------------------------------------------------------------
$ws.Range("H13").Select() | Out-Null

try { $wb.Styles("Standard").Name = "Normal" } catch {}
